# Refresh the "cryptos" price table (Sheet1, columns B:E, rows 2-51) with the
# latest scrape: updated Price (D) / Volume(1h) (E) figures for (almost)
# every coin, plus two rank swaps where a coin's rounded price moved past its
# neighbour: row 13/14 (Polkadot <-> WrappedEther) and row 50/51
# (EnergySwap <-> USDD).
#
# Column D holds prices formatted as plain text (e.g. "25.956.38",
# "0.0618", "1.01") rather than numbers - that's how the source sheet stores
# them (t="inlineStr"/shared-string cells, General number format). Excel's
# COM layer auto-coerces a purely-numeric-looking string typed into
# .Value into a real Number, which would silently drop things like the
# trailing zero in "0.250" or turn "1.01" into 1.01 (losing its text type).
# Set-TextValue works around that the way a human typist would: prefix the
# literal value with an apostrophe (forces Excel's "Text" quote-prefix
# entry mode) and then reset the cell's style back to Normal afterwards so
# no stray formatting is left behind - only the text value itself sticks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $cell = $ws.Range($range)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Numeric-looking text (e.g. price column D) - force text entry.
        $cell.Value = "'" + $value
        $cell.Style = 'Normal'
    } else {
        # Already unambiguous as text (URLs, names, padded "  +x.xx%  ").
        $cell.Value = $value
    }
}

# --- Row 2: Bitcoin ---------------------------------------------------
Set-TextValue 'D2' '25.956.38'
Set-TextValue 'E2' '  -0.27%  '

# --- Row 3: Ethereum ----------------------------------------------------
Set-TextValue 'D3' '1.627.69'
Set-TextValue 'E3' '  -0.98%  '

# --- Row 4: TetherUSD -----------------------------------------------------
Set-TextValue 'E4' '  -0.10%  '

# --- Row 5: BNB -----------------------------------------------------------
Set-TextValue 'D5' '214.19'
Set-TextValue 'E5' '  -0.82%  '

# --- Row 6: XRP -------------------------------------------------------------
Set-TextValue 'E6' '  -0.75%  '

# --- Row 7: USDC ------------------------------------------------------------
Set-TextValue 'E7' '  -0.01%  '

# --- Row 8: Cardano -----------------------------------------------------
Set-TextValue 'D8' '0.250'
Set-TextValue 'E8' '  -2.02%  '

# --- Row 9: Dogecoin ----------------------------------------------------
Set-TextValue 'D9' '0.0618'
Set-TextValue 'E9' '  -3.27%  '

# --- Row 10: Solana -----------------------------------------------------
Set-TextValue 'D10' '18.45'
Set-TextValue 'E10' '  -5.83%  '

# --- Row 11: TRON -----------------------------------------------------------
Set-TextValue 'E11' '  -1.22%  '

# --- Row 12: WrappedliquidstakedEther2.0 -----------------------------------
Set-TextValue 'D12' '1.854.57'
Set-TextValue 'E12' '  -0.89%  '

# --- Rows 13/14: Polkadot and WrappedEther swap rank -------------------
# WrappedEther now ranks above Polkadot, so row 13 becomes WrappedEther
# and row 14 becomes Polkadot (their Link/Price/Volume move along with them).
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.642.58'
Set-TextValue 'E13' '  -0.80%  '

Set-TextValue 'B14' 'Polkadot'
Set-TextValue 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D14' '4.18'
Set-TextValue 'E14' '  -2.10%  '

# --- Row 15: Polygon ------------------------------------------------------
Set-TextValue 'D15' '0.527'
Set-TextValue 'E15' '  -3.26%  '

# --- Row 16: WrappedBTC ---------------------------------------------------
Set-TextValue 'D16' '25.959.17'
Set-TextValue 'E16' '  -0.35%  '

# --- Row 17: ShibaInu -----------------------------------------------------
Set-TextValue 'D17' '0.0₃0739'
Set-TextValue 'E17' '  -3.06%  '

# --- Row 18: Litecoin -----------------------------------------------------
Set-TextValue 'D18' '61.31'
Set-TextValue 'E18' '  -3.34%  '

# --- Row 19: Dai ------------------------------------------------------------
Set-TextValue 'E19' '  -0.06%  '

# --- Row 20: BitcoinCash ----------------------------------------------------
Set-TextValue 'D20' '192.44'
Set-TextValue 'E20' '  -0.95%  '

# --- Row 21: Uniswap ------------------------------------------------------
Set-TextValue 'D21' '4.24'
Set-TextValue 'E21' '  -2.62%  '

# --- Row 22: Avalanche ----------------------------------------------------
Set-TextValue 'D22' '9.58'
Set-TextValue 'E22' '  -3.50%  '

# --- Row 23: Chainlink ----------------------------------------------------
Set-TextValue 'E23' '  -2.01%  '

# --- Row 24: Stellar ------------------------------------------------------
Set-TextValue 'E24' '  +0.56%  '

# --- Row 25: Monero -------------------------------------------------------
Set-TextValue 'D25' '143.81'
Set-TextValue 'E25' '  +0.40%  '

# --- Row 26: BinanceUSD -----------------------------------------------------
Set-TextValue 'E26' '  +0.04%  '

# --- Row 27: Toncoin ------------------------------------------------------
Set-TextValue 'E27' '  -3.72%  '

# --- Row 28: Cosmos -------------------------------------------------------
Set-TextValue 'D28' '6.73'
Set-TextValue 'E28' '  -2.16%  '

# --- Row 29: EthereumClassic ------------------------------------------------
Set-TextValue 'D29' '15.20'
Set-TextValue 'E29' '  -1.95%  '

# --- Row 30: PancakeSwap ----------------------------------------------------
Set-TextValue 'E30' '  -1.39%  '

# --- Row 31: Hedera -------------------------------------------------------
Set-TextValue 'E31' '  -2.09%  '

# --- Row 32: Filecoin -----------------------------------------------------
Set-TextValue 'E32' '  -4.10%  '

# --- Row 33: InternetComputer(DFINITY) -------------------------------------
Set-TextValue 'D33' '3.11'
Set-TextValue 'E33' '  -5.50%  '

# --- Row 34: LidoDAOToken ----------------------------------------------------
Set-TextValue 'E34' '  -2.80%  '

# --- Row 35: HuobiToken ---------------------------------------------------
Set-TextValue 'E35' '  -2.60%  '

# --- Row 36: Maker --------------------------------------------------------
Set-TextValue 'D36' '1.126.00'
Set-TextValue 'E36' '  -0.46%  '

# --- Row 37: ARBITRUM -----------------------------------------------------
Set-TextValue 'D37' '0.853'
Set-TextValue 'E37' '  -5.77%  '

# --- Row 38: MXToken ------------------------------------------------------
Set-TextValue 'E38' '  -1.56%  '

# --- Row 39: ImmutableX ---------------------------------------------------
Set-TextValue 'D39' '0.522'
Set-TextValue 'E39' '  -3.45%  '

# --- Row 40: VeChain ------------------------------------------------------
Set-TextValue 'E40' '  -2.33%  '

# --- Row 41: Quant --------------------------------------------------------
Set-TextValue 'D41' '98.17'
Set-TextValue 'E41' '  -1.03%  '

# --- Row 42: RocketPoolETH -------------------------------------------------
Set-TextValue 'D42' '1.765.80'
Set-TextValue 'E42' '  -0.85%  '

# --- Row 43: TrustWalletToken -----------------------------------------------
Set-TextValue 'E43' '  -4.28%  '

# --- Row 44: FraxShare ----------------------------------------------------
Set-TextValue 'D44' '5.17'
Set-TextValue 'E44' '  -5.29%  '

# --- Row 45: Cronos -------------------------------------------------------
Set-TextValue 'E45' '  +1.98%  '

# --- Row 46: Aave ---------------------------------------------------------
Set-TextValue 'D46' '54.40'

# --- Row 47: RenderToken ---------------------------------------------------
Set-TextValue 'D47' '1.47'
Set-TextValue 'E47' '  -1.10%  '

# --- Row 48: BabyDogeCoin --------------------------------------------------
Set-TextValue 'D48' '0.0₇0978'
Set-TextValue 'E48' '  -16.63%  '

# --- Row 49: Mantle -------------------------------------------------------
Set-TextValue 'E49' '  -0.57%  '

# --- Rows 50/51: EnergySwap and USDD swap rank --------------------------
# USDD now ranks above EnergySwap, so row 50 becomes USDD and row 51
# becomes EnergySwap (their Link/Price/Volume move along with them).
Set-TextValue 'B50' 'USDD'
Set-TextValue 'C50' 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue 'D50' '1.01'
Set-TextValue 'E50' '  +0.25%  '

Set-TextValue 'B51' 'EnergySwap'
Set-TextValue 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '7.47'
Set-TextValue 'E51' '  -3.70%  '
